$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "always have a page with @facs may or maynot have @subtype"
$ws.Range("D5").Value = "always have at least one div group in a page"
$ws.Range("D10").Value = "list must have a head if parent div group doesn't"

$ws.Range("A13").Select()
